$wb = $excel.ActiveWorkbook

# --- Rename existing sheet to "demand" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "demand"

# --- Update header row (row 1) ---
$ws1.Range("A1").Value = "hours"
$ws1.Range("B1").Value = "H2"
$ws1.Range("C1").Value = "Industrial heat"
$ws1.Range("D1").Value = "total"

# --- Copy header formatting (s="1") down column A/D to match new layout ---

$ws1.Range("A2").Value = "2019-01-01_00:00:00"
$ws1.Range("B2").Value = 14494062.7969296
$ws1.Range("C2").Value = 891387.0095024927
$ws1.Range("D2").Value = 15385449.8064321
$ws1.Range("A3").Value = "2019-01-31_10:00:00"
$ws1.Range("B3").Value = 10087332.83098327
$ws1.Range("C3").Value = 1308460.386135032
$ws1.Range("D3").Value = 11395793.2171183
$ws1.Range("A4").Value = "2019-03-02_20:00:00"
$ws1.Range("B4").Value = 10087512.49240243
$ws1.Range("C4").Value = 2979408.461069473
$ws1.Range("D4").Value = 13066920.9534719
$ws1.Range("A5").Value = "2019-04-02_06:00:00"
$ws1.Range("B5").Value = 6726528.25539057
$ws1.Range("C5").Value = 3249401.053300465
$ws1.Range("D5").Value = 9975929.308691034
$ws1.Range("A6").Value = "2019-05-02_16:00:00"
$ws1.Range("B6").Value = 11113028.48583499
$ws1.Range("C6").Value = 3358679.843611017
$ws1.Range("D6").Value = 14471708.32944601
$ws1.Range("A7").Value = "2019-06-02_02:00:00"
$ws1.Range("B7").Value = 8954696.288220024
$ws1.Range("C7").Value = 3349416.952809414
$ws1.Range("D7").Value = 12304113.24102944
$ws1.Range("A8").Value = "2019-07-02_12:00:00"
$ws1.Range("B8").Value = 9557711.991646467
$ws1.Range("C8").Value = 3340606.379506555
$ws1.Range("D8").Value = 12898318.37115302
$ws1.Range("A9").Value = "2019-08-01_22:00:00"
$ws1.Range("B9").Value = 7307068.730919966
$ws1.Range("C9").Value = 3355562.656824089
$ws1.Range("D9").Value = 10662631.38774406
$ws1.Range("A10").Value = "2019-09-01_08:00:00"
$ws1.Range("B10").Value = 8366500.548522449
$ws1.Range("C10").Value = 3358679.82090106
$ws1.Range("D10").Value = 11725180.36942351
$ws1.Range("A11").Value = "2019-10-01_18:00:00"
$ws1.Range("B11").Value = 4405940.086694176
$ws1.Range("C11").Value = 3028760.585553382
$ws1.Range("D11").Value = 7434700.672247559
$ws1.Range("A12").Value = "2019-11-01_04:00:00"
$ws1.Range("B12").Value = 13648237.51974842
$ws1.Range("C12").Value = 3285022.565868652
$ws1.Range("D12").Value = 16933260.08561707
$ws1.Range("A13").Value = "2019-12-01_14:00:00"
$ws1.Range("B13").Value = 15288168.45128457
$ws1.Range("C13").Value = 3275855.840447044
$ws1.Range("D13").Value = 18564024.29173161

# --- Apply header style (bold, border, centered) to A2:A13 and D1:D13 ---
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("A2:A13").PasteSpecial(-4122) | Out-Null
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Add new "prices" sheet after "demand" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "prices"
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Populate prices sheet ---
$ws2.Range("A1").Value = "hours"
$ws2.Range("B1").Value = "NED"
$ws2.Range("A2").Value = "2019-01-01_00:00:00"
$ws2.Range("B2").Value = 188.468197643773
$ws2.Range("A3").Value = "2019-01-31_10:00:00"
$ws2.Range("B3").Value = 160.2863019909899
$ws2.Range("A4").Value = "2019-03-02_20:00:00"
$ws2.Range("B4").Value = 110.7809744493958
$ws2.Range("A5").Value = "2019-04-02_06:00:00"
$ws2.Range("B5").Value = 55.79715927674411
$ws2.Range("A6").Value = "2019-05-02_16:00:00"
$ws2.Range("B6").Value = 55.4921684300246
$ws2.Range("A7").Value = "2019-06-02_02:00:00"
$ws2.Range("B7").Value = 55.6175675413002
$ws2.Range("A8").Value = "2019-07-02_12:00:00"
$ws2.Range("B8").Value = 55.5860344634459
$ws2.Range("A9").Value = "2019-08-01_22:00:00"
$ws2.Range("B9").Value = 55.75963361921727
$ws2.Range("A10").Value = "2019-09-01_08:00:00"
$ws2.Range("B10").Value = 63.5734283013198
$ws2.Range("A11").Value = "2019-10-01_18:00:00"
$ws2.Range("B11").Value = 73.40910145377882
$ws2.Range("A12").Value = "2019-11-01_04:00:00"
$ws2.Range("B12").Value = 70.42053777034438
$ws2.Range("A13").Value = "2019-12-01_14:00:00"
$ws2.Range("B13").Value = 46.07482878091285

# --- Apply header style to prices sheet A1:B1 and A2:A13 ---
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A2:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wb.Worksheets.Item(1).Activate()
